# Weekly crypto price/volume refresh (GitHub Actions scraper run, 2023-12-07 09:23:47 UTC).
# Source sheet stores every data cell as literal text (values like '43.391.69' or '0.0₃0961'
# are not valid Excel numbers, and the '%' columns keep their padding spaces), so plain numeric-
# looking strings are forced to Text via a leading apostrophe and then ClearFormats() strips the
# resulting @ number-format style back off, leaving the cell on the sheet's default (unstyled) xf
# -- exactly like the rest of the data cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.391.69'
$ws.Range("E2").Value = '  -0.90%  '

# Row 3
$ws.Range("D3").Value = '2.238.17'
$ws.Range("E3").Value = '  -1.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = "'230.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.39%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = "'0.641"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.29%  '

# Row 7
$ws.Range("D7").Value = "'63.05"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.99%  '

# Row 8
$ws.Range("E8").Value = '  +0.10%  '

# Row 9
$ws.Range("D9").Value = "'0.440"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.75%  '

# Row 10
$ws.Range("D10").Value = "'0.0957"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.40%  '

# Row 11
$ws.Range("D11").Value = "'56.84"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.39%  '

# Row 12
$ws.Range("D12").Value = "'26.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.94%  '

# Row 13
$ws.Range("D13").Value = "'0.104"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.55%  '

# Row 14
$ws.Range("D14").Value = '2.571.29'
$ws.Range("E14").Value = '  -1.77%  '

# Row 15
$ws.Range("D15").Value = "'15.41"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.01%  '

# Row 16
$ws.Range("D16").Value = "'6.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.63%  '

# Row 17
$ws.Range("D17").Value = "'0.826"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.95%  '

# Row 18
$ws.Range("D18").Value = '2.237.35'
$ws.Range("E18").Value = '  -2.07%  '

# Row 19
$ws.Range("D19").Value = '43.272.18'
$ws.Range("E19").Value = '  -0.99%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -3.34%  '

# Row 21
$ws.Range("D21").Value = "'72.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.11%  '

# Row 22
$ws.Range("D22").Value = "'6.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.70%  '

# Row 23
$ws.Range("D23").Value = "'245.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -7.70%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").Value = "'3.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +35.13%  '

# Row 26
$ws.Range("D26").Value = "'2.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.66%  '

# Row 27
$ws.Range("D27").Value = "'2.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.87%  '

# Row 28
$ws.Range("D28").Value = "'9.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.91%  '

# Row 29
$ws.Range("D29").Value = "'171.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.63%  '

# Row 30
$ws.Range("D30").Value = "'21.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.20%  '

# Row 31
$ws.Range("D31").Value = "'0.131"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.55%  '

# Row 32
$ws.Range("E32").Value = '  -3.46%  '

# Row 33
$ws.Range("D33").Value = "'0.125"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.42%  '

# Row 34
$ws.Range("D34").Value = "'4.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.98%  '

# Row 35
$ws.Range("D35").Value = "'0.0672"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.54%  '

# Row 36
$ws.Range("D36").Value = "'4.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.89%  '

# Row 37
$ws.Range("D37").Value = "'3.60"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.60%  '

# Row 38
$ws.Range("D38").Value = "'6.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -7.53%  '

# Row 39
$ws.Range("D39").Value = "'2.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.56%  '

# Row 40
$ws.Range("D40").Value = "'0.0250"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("D42").Value = "'8.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.75%  '

# Row 43
$ws.Range("D43").Value = "'4.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.40%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = "'97.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.74%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = "'16.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.46%  '

# Row 46
$ws.Range("D46").Value = "'1.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.30%  '

# Row 47
$ws.Range("D47").Value = "'0.0931"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.10%  '

# Row 48
$ws.Range("D48").Value = '1.445.82'
$ws.Range("E48").Value = '  -2.49%  '

# Row 49
$ws.Range("D49").Value = "'0.000207"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.45%  '

# Row 50
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = "'2.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.75%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = "'2.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.60%  '
